$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value would otherwise be auto-parsed as a
# number by Excel's type inference; force them to remain plain text so
# they keep their exact original string representation (e.g. '625.04',
# not 625.03999999999996), then drop back to the Normal style so no
# stray number-format is left behind on the cell.
$textCells = @('D5', 'D6', 'D11', 'D14', 'D15', 'D17', 'D19', 'D20', 'D21', 'D23', 'D24', 'D25', 'D26', 'D27', 'D30', 'D31', 'D34', 'D35', 'D37', 'D38', 'D40', 'D41', 'D43', 'D44', 'D45', 'D47', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '97.576.74'
$ws.Range('E2').Value = '  +3.29%  '
$ws.Range('D3').Value = '3.349.51'
$ws.Range('E3').Value = '  +7.65%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '246.12'
$ws.Range('E5').Value = '  +3.22%  '
$ws.Range('D6').Value = '625.04'
$ws.Range('E6').Value = '  +1.11%  '
$ws.Range('E7').Value = '  +0.72%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').Value = '3.349.94'
$ws.Range('E10').Value = '  +7.75%  '
$ws.Range('D11').Value = '0.799'
$ws.Range('E11').Value = '  -3.89%  '
$ws.Range('E12').Value = '  +1.29%  '
$ws.Range('D13').Value = '97.355.99'
$ws.Range('E13').Value = '  +3.49%  '
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').Value = '0.0000249'
$ws.Range('E14').Value = '  +1.80%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').Value = '35.67'
$ws.Range('E15').Value = '  +1.50%  '
$ws.Range('D16').Value = '3.963.66'
$ws.Range('E16').Value = '  +7.52%  '
$ws.Range('D17').Value = '5.53'
$ws.Range('E17').Value = '  +1.97%  '
$ws.Range('D18').Value = '3.356.64'
$ws.Range('E18').Value = '  +8.64%  '
$ws.Range('D19').Value = '3.63'
$ws.Range('E19').Value = '  -0.85%  '
$ws.Range('D20').Value = '15.25'
$ws.Range('E20').Value = '  +1.00%  '
$ws.Range('D21').Value = '493.93'
$ws.Range('E21').Value = '  +10.14%  '
$ws.Range('E22').Value = '  +5.71%  '
$ws.Range('D23').Value = '5.92'
$ws.Range('E23').Value = '  -2.00%  '
$ws.Range('D24').Value = '9.34'
$ws.Range('E24').Value = '  +3.18%  '
$ws.Range('D25').Value = '5.73'
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('D26').Value = '88.63'
$ws.Range('E26').Value = '  +2.67%  '
$ws.Range('D27').Value = '12.17'
$ws.Range('E27').Value = '  -2.26%  '
$ws.Range('D28').Value = '3.504.18'
$ws.Range('E28').Value = '  +7.10%  '
$ws.Range('E29').Value = '  +0.26%  '
$ws.Range('D30').Value = '0.184'
$ws.Range('E30').Value = '  +1.55%  '
$ws.Range('D31').Value = '0.241'
$ws.Range('E31').Value = '  -3.25%  '
$ws.Range('E32').Value = '  -3.62%  '
$ws.Range('E33').Value = '  -0.22%  '
$ws.Range('D34').Value = '9.41'
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('D35').Value = '27.76'
$ws.Range('E35').Value = '  +5.91%  '
$ws.Range('E36').Value = '  -4.20%  '
$ws.Range('D37').Value = '7.52'
$ws.Range('E37').Value = '  -5.22%  '
$ws.Range('D38').Value = '505.34'
$ws.Range('E38').Value = '  +4.95%  '
$ws.Range('E39').Value = '  +1.95%  '
$ws.Range('D40').Value = '24.74'
$ws.Range('E40').Value = '  +3.14%  '
$ws.Range('D41').Value = '0.456'
$ws.Range('E41').Value = '  -1.17%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').Value = '0.805'
$ws.Range('E43').Value = '  +16.11%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = '3.29'
$ws.Range('E44').Value = '  +1.50%  '
$ws.Range('D45').Value = '3.50'
$ws.Range('E45').Value = '  -7.22%  '
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').Value = '160.75'
$ws.Range('E47').Value = '  +0.10%  '
$ws.Range('E48').Value = '  +5.38%  '
$ws.Range('D49').Value = '4.59'
$ws.Range('E49').Value = '  +3.12%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = '0.0335'
$ws.Range('E50').Value = '  +3.60%  '
$ws.Range('B51').Value = 'ImmutableX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D51').Value = '1.36'
$ws.Range('E51').Value = '  +2.33%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
